$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.525.29'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.469.70'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = "'314.79"
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = "'92.00"
$ws.Range('E6').Value = '  -3.10%  '
$ws.Range('D7').Value = "'0.548"
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = "'0.514"
$ws.Range('E9').Value = '  +2.65%  '
$ws.Range('D10').Value = "'32.18"
$ws.Range('E10').Value = '  -4.40%  '
$ws.Range('D11').Value = "'0.0789"
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '2.850.76'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('D14').Value = "'6.83"
$ws.Range('E14').Value = '  -2.61%  '
$ws.Range('D15').Value = "'15.95"
$ws.Range('E15').Value = '  +3.16%  '
$ws.Range('D16').Value = '2.472.35'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('E17').Value = '  -3.06%  '
$ws.Range('D18').Value = '41.506.95'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  +1.81%  '
$ws.Range('D21').Value = "'71.27"
$ws.Range('E21').Value = '  +3.20%  '
$ws.Range('D22').Value = "'11.07"
$ws.Range('E22').Value = '  -2.37%  '
$ws.Range('D23').Value = "'235.80"
$ws.Range('E23').Value = '  -0.69%  '
$ws.Range('E24').Value = '  -1.69%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = "'1.89"
$ws.Range('E26').Value = '  -1.26%  '
$ws.Range('E27').Value = '  +1.41%  '
$ws.Range('E28').Value = '  -0.84%  '
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('D30').Value = "'35.23"
$ws.Range('E30').Value = '  -3.73%  '
$ws.Range('D31').Value = "'155.72"
$ws.Range('E31').Value = '  +2.26%  '
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').Value = "'17.24"
$ws.Range('E35').Value = '  -4.46%  '
$ws.Range('E36').Value = '  -7.33%  '
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('E38').Value = '  -0.81%  '
$ws.Range('D39').Value = "'1.77"
$ws.Range('E39').Value = '  -5.87%  '
$ws.Range('E40').Value = '  -12.61%  '
$ws.Range('D41').Value = "'4.03"
$ws.Range('E41').Value = '  -4.06%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '1.940.22'
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').Value = "'18.42"
$ws.Range('E45').Value = '  -7.23%  '
$ws.Range('D46').Value = "'2.92"
$ws.Range('E46').Value = '  -4.07%  '
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('D48').Value = '2.708.50'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').Value = "'96.80"
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').Value = "'66.94"
$ws.Range('E50').Value = '  -4.61%  '
$ws.Range('D51').Value = "'52.69"
$ws.Range('E51').Value = '  +3.00%  '

Write-Output "Updated cryptos list"